# Update "想去人数" (interested count) values for a handful of events.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) rows 3-6, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 82
$wsExpo.Range("F4").Value = 2353
$wsExpo.Range("F5").Value = 25
$wsExpo.Range("F6").Value = 517

# Sheet "全部类型" (all types) rows 5-8, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 82
$wsAll.Range("F6").Value = 2353
$wsAll.Range("F7").Value = 25
$wsAll.Range("F8").Value = 517
